$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback Datetime for first data row (row 2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-12 16:52:30"
$wsZh.Range("H2").Value = "2016-03-12 16:52:46"

# de-de sheet: update Correspond Handoff/Handback Datetime for first data row (row 2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-12 16:52:33"
$wsDe.Range("H2").Value = "2016-03-12 16:52:51"
